# #91: fixes serialization issue
#
# Slide 1 changes:
#  - "Rectangle 22": grow its height and add the missing
#    "MinioArchiverProcessor" label.
#  - "Connector: Curved 32" (glued to Rectangle 22's bottom connection
#    site): re-anchor/re-size so it still meets the rectangle's new,
#    taller edge.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Rectangle 22 -------------------------------------------------------
$rect = $s.Shapes.Item("Rectangle 22")

# Grow the box's height; left/top/width are unchanged.
$rect.Height = 196.94582677165354

# Add the missing label text. Use InsertBefore (rather than assigning
# .Text) so the paragraph's existing endParaRPr is preserved instead of
# being replaced outright.
$tr = $rect.TextFrame.TextRange
$run = $tr.InsertBefore("MinioArchiverProcessor")
$run.LanguageID = "en-GB"

# --- Connector: Curved 32 ------------------------------------------------
$conn = $s.Shapes.Item("Connector: Curved 32")

# Keep it glued to Rectangle 22's (now lower) connection point: shift the
# top down and shrink the height; left/width are unchanged.
$conn.Top = 392.391983503937
$conn.Height = 59.45653743307087
